# Identity-AuthorizeExerciseMVC.pptx
# Merge split runs back into single runs (text content unchanged,
# only the run/formatting boundaries collapse to match the author's
# final edit).

$p = $ppt.ActivePresentation

# --- Slide 2: title paragraph "Authorization in MVC 5 using ASP.NET Identity" ---
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$full2 = $tr2.Text
$target2 = "Authorization in MVC 5 using ASP.NET Identity"
$idx2 = $full2.IndexOf($target2) + 1
$sub2 = $tr2.Characters($idx2, $target2.Length)
$sub2.Text = $target2

# --- Slide 4: body paragraph about the authorization attribute ---
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange
$full4 = $tr4.Text
$target4 = "The authorization attribute can also be used to check if an identity is in a specific role or is a user "
$idx4 = $full4.IndexOf($target4) + 1
$sub4 = $tr4.Characters($idx4, $target4.Length)
$sub4.Text = $target4
